# Applies a 3-way rotation of data among rows 2, 3 and 4 for columns
# A, B, E, F, G, H, Q, R:
#   new row2 = old row4
#   new row3 = old row2
#   new row4 = old row3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")
$rows = @(2, 3, 4)

# Capture the "old" values for the relevant columns/rows before overwriting
$old = @{}
foreach ($r in $rows) {
    $old[$r] = @{}
    foreach ($col in $cols) {
        $old[$r][$col] = $ws.Range("$col$r").Value2
    }
}

# Mapping: new row <- old row
$rotation = @{ 2 = 4; 3 = 2; 4 = 3 }

foreach ($r in $rows) {
    $srcRow = $rotation[$r]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $old[$srcRow][$col]
    }
}
